# This script reshuffles the "trial" rows (rows 2-41) of the kitchens
# categorization input file. Columns G,H,I,K,L,M,N,O,P,Q,R,S,T,U,V of each
# destination row are replaced with the corresponding values taken from a
# source row (a permutation of rows 2-41, snapshotted BEFORE any writes so
# that the copy is based entirely on the original data). Column F
# (trial_total) is simply reset to the row's position within the block
# (1..40), matching column E (trial_block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Destination row -> source row mapping (permutation of rows 2..41)
$rowMap = @{
    2=35; 3=36; 4=8; 5=25; 6=29; 7=3; 8=13; 9=26; 10=30;
    11=21; 12=10; 13=7; 14=15; 15=18; 16=2; 17=32; 18=4; 19=38; 20=11;
    21=39; 22=24; 23=23; 24=41; 25=28; 26=9; 27=40; 28=27; 29=5; 30=33;
    31=12; 32=19; 33=6; 34=20; 35=14; 36=37; 37=34; 38=17; 39=16; 40=22;
    41=31
}

# Columns (by number) whose contents get copied from the source row.
# G=7 H=8 I=9 K=11 L=12 M=13 N=14 O=15 P=16 Q=17 R=18 S=19 T=20 U=21 V=22
$copyCols = @(7, 8, 9, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22)

# --- Pass 1: snapshot all the "before" values for rows 2..41 ---
$snapshot = @{}
for ($r = 2; $r -le 41; $r++) {
    $rowVals = @{}
    foreach ($c in $copyCols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# --- Pass 2: write the new values using the snapshotted source rows ---
for ($r = 2; $r -le 41; $r++) {
    $src = $rowMap[$r]
    $srcVals = $snapshot[$src]
    foreach ($c in $copyCols) {
        $ws.Cells.Item($r, $c).Value = $srcVals[$c]
    }
    # trial_total (F) becomes the 1-based position within the block
    $ws.Cells.Item($r, 6).Value = ($r - 1)
}

Write-Host "Row reshuffle complete"
